$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data point to the end of the Ticker list in column A
$ws.Range("A65").Value = "GRT-USD"
